# Updates cryptocurrency price/volume data in the "cryptos" worksheet
# per the GitHub Actions scheduled data refresh (Tue Jan 24 23:41:39 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) / Volume(1h) (column E) values, keyed by row number.
$updates = @{
    2 = @{ 'D' = '302.20'; 'E' = '-1.20%' }
    3 = @{ 'D' = '35.36'; 'E' = '-2.27%' }
    4 = @{ 'D' = '5.008'; 'E' = '-1.19%' }
    5 = @{ 'D' = '0.07838'; 'E' = '-1.61%' }
    6 = @{ 'D' = '1.829'; 'E' = '-15.62%' }
    7 = @{ 'E' = '-1.88%' }
    8 = @{ 'D' = '7.812'; 'E' = '-2.71%' }
    9 = @{ 'D' = '2.840'; 'E' = '7.92%' }
    10 = @{ 'D' = '0.9222'; 'E' = '-0.85%' }
    11 = @{ 'D' = '0.1079'; 'E' = '9.30%' }
    12 = @{ 'D' = '0.1859'; 'E' = '-0.36%' }
    13 = @{ 'D' = '0.09364'; 'E' = '3.67%' }
    14 = @{ 'D' = '0.03579'; 'E' = '-1.15%' }
    15 = @{ 'D' = '0.09937'; 'E' = '0.01%' }
    16 = @{ 'D' = '0.001404'; 'E' = '-2.37%' }
    17 = @{ 'D' = '0.005713'; 'E' = '0.64%' }
    18 = @{ 'D' = '3.461'; 'E' = '0.30%' }
    19 = @{ 'D' = '0.3430'; 'E' = '1.74%' }
    20 = @{ 'E' = '-4.57%' }
    21 = @{ 'D' = '5.130'; 'E' = '0.80%' }
    23 = @{ 'D' = '0.04568'; 'E' = '-0.25%' }
    24 = @{ 'E' = '-0.93%' }
    25 = @{ 'D' = '0.004652'; 'E' = '-2.20%' }
    26 = @{ 'E' = '-3.52%' }
    27 = @{ 'D' = '0.0004470'; 'E' = '-5.66%' }
    39 = @{ 'D' = '0.01890'; 'E' = '-2.72%' }
    40 = @{ 'D' = '0.04711'; 'E' = '-3.89%' }
    41 = @{ 'D' = '0.007574'; 'E' = '-3.10%' }
    42 = @{ 'D' = '0.01001'; 'E' = '28.10%' }
    43 = @{ 'D' = '0.1333'; 'E' = '-4.28%' }
    44 = @{ 'E' = '1.18%' }
    45 = @{ 'E' = '-1.52%' }
    46 = @{ 'D' = '0.00006277'; 'E' = '0.83%' }
    47 = @{ 'E' = '0.25%' }
    48 = @{ 'E' = '27.35%' }
    49 = @{ 'D' = '0.001306'; 'E' = '-27.59%' }
    50 = @{ 'D' = '0.00002110'; 'E' = '0.25%' }
    51 = @{ 'D' = '0.0002010'; 'E' = '0.25%' }
}

# Column D/E cells store plain text (e.g. "302.20", "-1.20%"), so force
# a text number format before assigning -- otherwise Excel would coerce
# the strings into numeric/percentage values and drop information such
# as trailing zeros (e.g. "302.20" -> 302.2).
foreach ($row in $updates.Keys) {
    $rowUpdates = $updates[$row]
    foreach ($col in $rowUpdates.Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $rowUpdates[$col]
    }
}
